$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 130  # ALC!H9
$ws.Cells.Item(9, 9).Value = 134.28572  # ALC!I9
$ws.Cells.Item(9, 11).Value = 134.28572  # ALC!K9
$ws.Cells.Item(9, 13).Value = 34.71428  # ALC!M9
$ws.Cells.Item(28, 8).Value = 634.96  # ALC!H28
$ws.Cells.Item(28, 9).Value = 387.33334  # ALC!I28
$ws.Cells.Item(28, 10).Value = 1271.7142  # ALC!J28
$ws.Cells.Item(28, 11).Value = 387.33334  # ALC!K28
$ws.Cells.Item(28, 12).Value = 1271.7142  # ALC!L28
$ws.Cells.Item(28, 13).Value = 97.66665999999998  # ALC!M28
$ws.Cells.Item(28, 14).Value = -2241.7142  # ALC!N28
$ws.Cells.Item(33, 8).Value = 6245811  # ALC!H33
$ws.Cells.Item(33, 9).Value = 11103500  # ALC!I33
$ws.Cells.Item(33, 11).Value = 11103500  # ALC!K33
$ws.Cells.Item(33, 13).Value = -11103271  # ALC!M33
$ws.Cells.Item(62, 8).Value = 3548.6924  # ALC!H62
$ws.Cells.Item(62, 9).Value = 2188.0833  # ALC!I62
$ws.Cells.Item(62, 11).Value = 2188.0833  # ALC!K62
$ws.Cells.Item(62, 13).Value = -1564.0833  # ALC!M62
$ws.Cells.Item(65, 8).Value = 3548.6924  # ALC!H65
$ws.Cells.Item(65, 9).Value = 2188.0833  # ALC!I65
$ws.Cells.Item(65, 11).Value = 10940.4165  # ALC!K65
$ws.Cells.Item(65, 13).Value = -7820.416499999999  # ALC!M65
$ws.Cells.Item(107, 8).Value = 567  # ALC!H107
$ws.Cells.Item(107, 9).Value = 567  # ALC!I107
$ws.Cells.Item(107, 10).Value = 0  # ALC!J107
$ws.Cells.Item(107, 11).Value = 567  # ALC!K107
$ws.Cells.Item(107, 12).Value = 0  # ALC!L107
$ws.Cells.Item(107, 13).Value = 1353  # ALC!M107
$ws.Cells.Item(107, 14).Value = $null  # ALC!N107
$ws.Cells.Item(111, 8).Value = 3468.5715  # ALC!H111
$ws.Cells.Item(111, 9).Value = 3160  # ALC!I111
$ws.Cells.Item(111, 11).Value = 9480  # ALC!K111
$ws.Cells.Item(111, 13).Value = -6413  # ALC!M111
$ws.Cells.Item(113, 8).Value = 2936.8333  # ALC!H113
$ws.Cells.Item(113, 9).Value = 2840  # ALC!I113
$ws.Cells.Item(113, 11).Value = 2840  # ALC!K113
$ws.Cells.Item(113, 13).Value = 414  # ALC!M113
$ws.Cells.Item(116, 8).Value = 2782.5386  # ALC!H116
$ws.Cells.Item(116, 10).Value = 2904.6  # ALC!J116
$ws.Cells.Item(116, 12).Value = 2904.6  # ALC!L116
$ws.Cells.Item(116, 14).Value = -9788.6  # ALC!N116
$ws.Cells.Item(129, 8).Value = 988.8269  # ALC!H129
$ws.Cells.Item(129, 10).Value = 1047.6809  # ALC!J129
$ws.Cells.Item(129, 12).Value = 3143.0427  # ALC!L129
$ws.Cells.Item(129, 14).Value = -13143.0427  # ALC!N129
$ws.Cells.Item(137, 8).Value = 879859.75  # ALC!H137
$ws.Cells.Item(137, 9).Value = 3170.3157  # ALC!I137
$ws.Cells.Item(137, 10).Value = 1756549.2  # ALC!J137
$ws.Cells.Item(137, 11).Value = 9510.947100000001  # ALC!K137
$ws.Cells.Item(137, 12).Value = 5269647.6  # ALC!L137
$ws.Cells.Item(137, 13).Value = -6960.947100000001  # ALC!M137
$ws.Cells.Item(137, 14).Value = -5274747.6  # ALC!N137
$ws.Cells.Item(139, 8).Value = 78000  # ALC!H139
$ws.Cells.Item(139, 10).Value = 78000  # ALC!J139
$ws.Cells.Item(139, 12).Value = 78000  # ALC!L139
$ws.Cells.Item(139, 14).Value = -88280  # ALC!N139
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 900  # ARM!H22
$ws.Cells.Item(22, 9).Value = 900  # ARM!I22
$ws.Cells.Item(22, 10).Value = 0  # ARM!J22
$ws.Cells.Item(22, 11).Value = 900  # ARM!K22
$ws.Cells.Item(22, 12).Value = 0  # ARM!L22
$ws.Cells.Item(22, 13).Value = -601  # ARM!M22
$ws.Cells.Item(22, 14).Value = $null  # ARM!N22
$ws.Cells.Item(32, 8).Value = 21110.203  # ARM!H32
$ws.Cells.Item(32, 9).Value = 26546.951  # ARM!I32
$ws.Cells.Item(32, 10).Value = 3963.5386  # ARM!J32
$ws.Cells.Item(32, 11).Value = 26546.951  # ARM!K32
$ws.Cells.Item(32, 12).Value = 3963.5386  # ARM!L32
$ws.Cells.Item(32, 13).Value = -26259.951  # ARM!M32
$ws.Cells.Item(32, 14).Value = -4537.5386  # ARM!N32
$ws.Cells.Item(74, 8).Value = 5534.7715  # ARM!H74
$ws.Cells.Item(74, 9).Value = 2390.6667  # ARM!I74
$ws.Cells.Item(74, 11).Value = 2390.6667  # ARM!K74
$ws.Cells.Item(74, 13).Value = -1516.6667  # ARM!M74
$ws.Cells.Item(77, 8).Value = 5534.7715  # ARM!H77
$ws.Cells.Item(77, 9).Value = 2390.6667  # ARM!I77
$ws.Cells.Item(77, 11).Value = 11953.3335  # ARM!K77
$ws.Cells.Item(77, 13).Value = -7585.333500000001  # ARM!M77
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(63, 8).Value = 38271  # BSM!H63
$ws.Cells.Item(63, 10).Value = 38271  # BSM!J63
$ws.Cells.Item(63, 12).Value = 38271  # BSM!L63
$ws.Cells.Item(63, 14).Value = -39643  # BSM!N63
$ws.Cells.Item(66, 8).Value = 38271  # BSM!H66
$ws.Cells.Item(66, 10).Value = 38271  # BSM!J66
$ws.Cells.Item(66, 12).Value = 114813  # BSM!L66
$ws.Cells.Item(66, 14).Value = -121677  # BSM!N66
$ws.Cells.Item(99, 8).Value = 1875.6471  # BSM!H99
$ws.Cells.Item(99, 9).Value = 1089.909  # BSM!I99
$ws.Cells.Item(99, 10).Value = 3316.1667  # BSM!J99
$ws.Cells.Item(99, 11).Value = 1089.909  # BSM!K99
$ws.Cells.Item(99, 12).Value = 3316.1667  # BSM!L99
$ws.Cells.Item(99, 13).Value = 408.0909999999999  # BSM!M99
$ws.Cells.Item(99, 14).Value = -6312.1667  # BSM!N99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 90  # CRP!H7
$ws.Cells.Item(7, 9).Value = 30  # CRP!I7
$ws.Cells.Item(7, 10).Value = 150  # CRP!J7
$ws.Cells.Item(7, 11).Value = 30  # CRP!K7
$ws.Cells.Item(7, 12).Value = 150  # CRP!L7
$ws.Cells.Item(7, 13).Value = 83  # CRP!M7
$ws.Cells.Item(7, 14).Value = -376  # CRP!N7
$ws.Cells.Item(31, 8).Value = 653026.3  # CRP!H31
$ws.Cells.Item(31, 9).Value = 5258.826  # CRP!I31
$ws.Cells.Item(31, 10).Value = 1091222  # CRP!J31
$ws.Cells.Item(31, 11).Value = 5258.826  # CRP!K31
$ws.Cells.Item(31, 12).Value = 1091222  # CRP!L31
$ws.Cells.Item(31, 13).Value = -4963.826  # CRP!M31
$ws.Cells.Item(31, 14).Value = -1091812  # CRP!N31
$ws.Cells.Item(34, 8).Value = 653026.3  # CRP!H34
$ws.Cells.Item(34, 9).Value = 5258.826  # CRP!I34
$ws.Cells.Item(34, 10).Value = 1091222  # CRP!J34
$ws.Cells.Item(34, 11).Value = 5258.826  # CRP!K34
$ws.Cells.Item(34, 12).Value = 1091222  # CRP!L34
$ws.Cells.Item(34, 13).Value = -5056.826  # CRP!M34
$ws.Cells.Item(34, 14).Value = -1091626  # CRP!N34
$ws.Cells.Item(105, 8).Value = 724.6667  # CRP!H105
$ws.Cells.Item(105, 9).Value = 724.6667  # CRP!I105
$ws.Cells.Item(105, 11).Value = 724.6667  # CRP!K105
$ws.Cells.Item(105, 13).Value = 1022.3333  # CRP!M105
$ws.Cells.Item(134, 8).Value = 2424.658  # CRP!H134
$ws.Cells.Item(134, 9).Value = 1918.1786  # CRP!I134
$ws.Cells.Item(134, 10).Value = 3842.8  # CRP!J134
$ws.Cells.Item(134, 11).Value = 5754.5358  # CRP!K134
$ws.Cells.Item(134, 12).Value = 11528.4  # CRP!L134
$ws.Cells.Item(134, 13).Value = -3219.5358  # CRP!M134
$ws.Cells.Item(134, 14).Value = -16598.4  # CRP!N134
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 2630.8167  # CUL!H68
$ws.Cells.Item(68, 10).Value = 5104.36  # CUL!J68
$ws.Cells.Item(68, 12).Value = 15313.08  # CUL!L68
$ws.Cells.Item(68, 14).Value = -16935.08  # CUL!N68
$ws.Cells.Item(70, 8).Value = 4700.3076  # CUL!H70
$ws.Cells.Item(70, 9).Value = 2204  # CUL!I70
$ws.Cells.Item(70, 10).Value = 5449.2  # CUL!J70
$ws.Cells.Item(70, 11).Value = 6612  # CUL!K70
$ws.Cells.Item(70, 12).Value = 16347.6  # CUL!L70
$ws.Cells.Item(70, 13).Value = -6297  # CUL!M70
$ws.Cells.Item(70, 14).Value = -16977.6  # CUL!N70
$ws.Cells.Item(71, 8).Value = 2630.8167  # CUL!H71
$ws.Cells.Item(71, 10).Value = 5104.36  # CUL!J71
$ws.Cells.Item(71, 12).Value = 45939.24  # CUL!L71
$ws.Cells.Item(71, 14).Value = -54051.24  # CUL!N71
$ws.Cells.Item(73, 8).Value = 4700.3076  # CUL!H73
$ws.Cells.Item(73, 9).Value = 2204  # CUL!I73
$ws.Cells.Item(73, 10).Value = 5449.2  # CUL!J73
$ws.Cells.Item(73, 11).Value = 6612  # CUL!K73
$ws.Cells.Item(73, 12).Value = 16347.6  # CUL!L73
$ws.Cells.Item(73, 13).Value = -5520  # CUL!M73
$ws.Cells.Item(73, 14).Value = -18531.6  # CUL!N73
$ws.Cells.Item(87, 8).Value = 5279.6665  # CUL!H87
$ws.Cells.Item(87, 9).Value = 4500.5  # CUL!I87
$ws.Cells.Item(87, 10).Value = 5502.2856  # CUL!J87
$ws.Cells.Item(87, 11).Value = 13501.5  # CUL!K87
$ws.Cells.Item(87, 12).Value = 16506.8568  # CUL!L87
$ws.Cells.Item(87, 13).Value = -12253.5  # CUL!M87
$ws.Cells.Item(87, 14).Value = -19002.8568  # CUL!N87
$ws.Cells.Item(90, 8).Value = 5279.6665  # CUL!H90
$ws.Cells.Item(90, 9).Value = 4500.5  # CUL!I90
$ws.Cells.Item(90, 10).Value = 5502.2856  # CUL!J90
$ws.Cells.Item(90, 11).Value = 40504.5  # CUL!K90
$ws.Cells.Item(90, 12).Value = 49520.5704  # CUL!L90
$ws.Cells.Item(90, 13).Value = -34264.5  # CUL!M90
$ws.Cells.Item(90, 14).Value = -62000.5704  # CUL!N90
$ws.Cells.Item(101, 8).Value = 6950.4443  # CUL!H101
$ws.Cells.Item(101, 10).Value = 6950.4443  # CUL!J101
$ws.Cells.Item(101, 12).Value = 20851.3329  # CUL!L101
$ws.Cells.Item(101, 14).Value = -25719.3329  # CUL!N101
$ws.Cells.Item(107, 8).Value = 681.1607  # CUL!H107
$ws.Cells.Item(107, 10).Value = 1786  # CUL!J107
$ws.Cells.Item(107, 12).Value = 5358  # CUL!L107
$ws.Cells.Item(107, 14).Value = -9198  # CUL!N107
$ws.Cells.Item(113, 8).Value = 370.27  # CUL!H113
$ws.Cells.Item(113, 9).Value = 409.25583  # CUL!I113
$ws.Cells.Item(113, 10).Value = 340.85965  # CUL!J113
$ws.Cells.Item(113, 11).Value = 1227.76749  # CUL!K113
$ws.Cells.Item(113, 12).Value = 1022.57895  # CUL!L113
$ws.Cells.Item(113, 13).Value = 942.23251  # CUL!M113
$ws.Cells.Item(113, 14).Value = -5362.57895  # CUL!N113
$ws.Cells.Item(114, 8).Value = 822.2727  # CUL!H114
$ws.Cells.Item(114, 9).Value = 419  # CUL!I114
$ws.Cells.Item(114, 10).Value = 1306.2  # CUL!J114
$ws.Cells.Item(114, 11).Value = 1257  # CUL!K114
$ws.Cells.Item(114, 12).Value = 3918.6  # CUL!L114
$ws.Cells.Item(114, 13).Value = 1997  # CUL!M114
$ws.Cells.Item(114, 14).Value = -10426.6  # CUL!N114
$ws.Cells.Item(131, 8).Value = 21008.125  # CUL!H131
$ws.Cells.Item(131, 10).Value = 23267.441  # CUL!J131
$ws.Cells.Item(131, 12).Value = 69802.323  # CUL!L131
$ws.Cells.Item(131, 14).Value = -79882.323  # CUL!N131
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 10000  # GSM!H2
$ws.Cells.Item(2, 9).Value = 0  # GSM!I2
$ws.Cells.Item(2, 11).Value = 0  # GSM!K2
$ws.Cells.Item(2, 13).Value = $null  # GSM!M2
$ws.Cells.Item(20, 8).Value = 10000000  # GSM!H20
$ws.Cells.Item(20, 10).Value = 0  # GSM!J20
$ws.Cells.Item(20, 12).Value = 0  # GSM!L20
$ws.Cells.Item(20, 14).Value = $null  # GSM!N20
$ws.Cells.Item(24, 8).Value = 13333333  # GSM!H24
$ws.Cells.Item(24, 9).Value = 13333333  # GSM!I24
$ws.Cells.Item(24, 11).Value = 13333333  # GSM!K24
$ws.Cells.Item(24, 13).Value = -13333160  # GSM!M24
$ws.Cells.Item(70, 8).Value = 5846.4443  # GSM!H70
$ws.Cells.Item(70, 9).Value = 5609.0713  # GSM!I70
$ws.Cells.Item(70, 11).Value = 5609.0713  # GSM!K70
$ws.Cells.Item(70, 13).Value = -5339.0713  # GSM!M70
$ws.Cells.Item(73, 8).Value = 5846.4443  # GSM!H73
$ws.Cells.Item(73, 9).Value = 5609.0713  # GSM!I73
$ws.Cells.Item(73, 11).Value = 5609.0713  # GSM!K73
$ws.Cells.Item(73, 13).Value = -4673.0713  # GSM!M73
$ws.Cells.Item(98, 8).Value = 0  # GSM!H98
$ws.Cells.Item(98, 10).Value = 0  # GSM!J98
$ws.Cells.Item(98, 12).Value = 0  # GSM!L98
$ws.Cells.Item(98, 14).Value = $null  # GSM!N98
$ws.Cells.Item(102, 8).Value = 3798.2727  # GSM!H102
$ws.Cells.Item(102, 9).Value = 3511.1667  # GSM!I102
$ws.Cells.Item(102, 11).Value = 3511.1667  # GSM!K102
$ws.Cells.Item(102, 13).Value = -1889.1667  # GSM!M102
$ws.Cells.Item(113, 8).Value = 1888.92  # GSM!H113
$ws.Cells.Item(113, 9).Value = 1781.875  # GSM!I113
$ws.Cells.Item(113, 10).Value = 2079.2222  # GSM!J113
$ws.Cells.Item(113, 11).Value = 1781.875  # GSM!K113
$ws.Cells.Item(113, 12).Value = 2079.2222  # GSM!L113
$ws.Cells.Item(113, 13).Value = 388.125  # GSM!M113
$ws.Cells.Item(113, 14).Value = -6419.2222  # GSM!N113
$ws.Cells.Item(132, 8).Value = 8118.3  # GSM!H132
$ws.Cells.Item(132, 9).Value = 3276.5715  # GSM!I132
$ws.Cells.Item(132, 10).Value = 10725.385  # GSM!J132
$ws.Cells.Item(132, 11).Value = 9829.7145  # GSM!K132
$ws.Cells.Item(132, 12).Value = 32176.155  # GSM!L132
$ws.Cells.Item(132, 13).Value = -7299.7145  # GSM!M132
$ws.Cells.Item(132, 14).Value = -37236.155  # GSM!N132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1291.2222  # LTW!H16
$ws.Cells.Item(16, 9).Value = 949.5  # LTW!I16
$ws.Cells.Item(16, 10).Value = 1388.8572  # LTW!J16
$ws.Cells.Item(16, 11).Value = 949.5  # LTW!K16
$ws.Cells.Item(16, 12).Value = 1388.8572  # LTW!L16
$ws.Cells.Item(16, 13).Value = -779.5  # LTW!M16
$ws.Cells.Item(16, 14).Value = -1728.8572  # LTW!N16
$ws.Cells.Item(40, 8).Value = 4401.8  # LTW!H40
$ws.Cells.Item(40, 9).Value = 3001  # LTW!I40
$ws.Cells.Item(40, 10).Value = 10005  # LTW!J40
$ws.Cells.Item(40, 11).Value = 3001  # LTW!K40
$ws.Cells.Item(40, 12).Value = 10005  # LTW!L40
$ws.Cells.Item(40, 13).Value = -2865  # LTW!M40
$ws.Cells.Item(40, 14).Value = -10277  # LTW!N40
$ws.Cells.Item(46, 8).Value = 857.7273  # LTW!H46
$ws.Cells.Item(46, 9).Value = 615  # LTW!I46
$ws.Cells.Item(46, 11).Value = 615  # LTW!K46
$ws.Cells.Item(46, 13).Value = -427  # LTW!M46
$ws.Cells.Item(120, 8).Value = 29999.666  # LTW!H120
$ws.Cells.Item(120, 10).Value = 29999.666  # LTW!J120
$ws.Cells.Item(120, 12).Value = 29999.666  # LTW!L120
$ws.Cells.Item(120, 14).Value = -39675.666  # LTW!N120
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1533.4445  # WVR!H107
$ws.Cells.Item(107, 9).Value = 991.25  # WVR!I107
$ws.Cells.Item(107, 10).Value = 1967.2  # WVR!J107
$ws.Cells.Item(107, 11).Value = 2973.75  # WVR!K107
$ws.Cells.Item(107, 12).Value = 5901.6  # WVR!L107
$ws.Cells.Item(107, 13).Value = -1053.75  # WVR!M107
$ws.Cells.Item(107, 14).Value = -9741.6  # WVR!N107
$ws.Cells.Item(122, 8).Value = 1960  # WVR!H122
$ws.Cells.Item(122, 9).Value = 1960  # WVR!I122
$ws.Cells.Item(122, 11).Value = 5880  # WVR!K122
$ws.Cells.Item(122, 13).Value = -3430  # WVR!M122
